$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds/score-line values for three fixture rows (27, 28, 44)
# on the weekly FlashScore odds sheet, per the upstream refresh.

# Row 27 updates
$ws.Range("G27").Value = 2.15
$ws.Range("H27").Value = 3.1
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 4.5
$ws.Range("M27").Value = 1.1
$ws.Range("N27").Value = 7
$ws.Range("O27").Value = 1.5
$ws.Range("P27").Value = 2.5
$ws.Range("Q27").Value = 2.6
$ws.Range("R27").Value = 1.48
$ws.Range("S27").Value = 1.57
$ws.Range("T27").Value = 2.25
$ws.Range("U27").Value = 2.2
$ws.Range("V27").Value = 1.62
$ws.Range("AA27").Value = 21
$ws.Range("AD27").Value = 6
$ws.Range("AE27").Value = 19
$ws.Range("AG27").Value = 8
$ws.Range("AH27").Value = 17
$ws.Range("AK27").Value = 41
$ws.Range("AP27").Value = 41
$ws.Range("AS27").Value = 2.25
$ws.Range("BB27").Value = 151

# Row 28 updates
$ws.Range("G28").Value = 3.7
$ws.Range("H28").Value = 3.5
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = 2.25
$ws.Range("L28").Value = 2.63
$ws.Range("M28").Value = 1.04
$ws.Range("N28").Value = 12
$ws.Range("O28").Value = 1.22
$ws.Range("P28").Value = 4
$ws.Range("Q28").Value = 1.8
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 1.36
$ws.Range("T28").Value = 3
$ws.Range("U28").Value = 1.67
$ws.Range("V28").Value = 2.1
$ws.Range("W28").Value = 12
$ws.Range("X28").Value = 19
$ws.Range("Y28").Value = 13
$ws.Range("Z28").Value = 41
$ws.Range("AC28").Value = 12
$ws.Range("AH28").Value = 10
$ws.Range("AI28").Value = 8.5
$ws.Range("AJ28").Value = 19
$ws.Range("AK28").Value = 15
$ws.Range("AL28").Value = 23
$ws.Range("AM28").Value = 5.5
$ws.Range("AN28").Value = 19
$ws.Range("AS28").Value = 3
$ws.Range("AV28").Value = 4
$ws.Range("AW28").Value = 11
$ws.Range("AX28").Value = 19
$ws.Range("AY28").Value = 34
$ws.Range("BA28").Value = 501
$ws.Range("BC28").Value = 126
$ws.Range("BD28").Value = 151

# Row 44 updates
$ws.Range("H44").Value = 3.75
$ws.Range("I44").Value = 3.7
$ws.Range("J44").Value = 2.5
$ws.Range("K44").Value = 2.25
$ws.Range("M44").Value = 1.04
$ws.Range("N44").Value = 13
$ws.Range("O44").Value = 1.25
$ws.Range("P44").Value = 3.75
$ws.Range("Q44").Value = 1.8
$ws.Range("R44").Value = 2
$ws.Range("S44").Value = 1.36
$ws.Range("T44").Value = 3
$ws.Range("U44").Value = 1.67
$ws.Range("V44").Value = 2.1
$ws.Range("W44").Value = 8
$ws.Range("X44").Value = 9.5
$ws.Range("AB44").Value = 23
$ws.Range("AC44").Value = 12
$ws.Range("AF44").Value = 41
$ws.Range("AG44").Value = 12
$ws.Range("AH44").Value = 21
$ws.Range("AL44").Value = 34
$ws.Range("AO44").Value = 19
$ws.Range("AR44").Value = 126
$ws.Range("AS44").Value = 3
$ws.Range("AT44").Value = 7.5
$ws.Range("AX44").Value = 26
$ws.Range("AZ44").Value = 81
$ws.Range("BC44").Value = 151
$ws.Range("BD44").Value = 201

Write-Host "Updated 89 cells across rows 27, 28, 44"
